$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eetu Pihamäki")

$ws.Range("A18").Value = 43383
$ws.Range("B18").Value = 0.5
$ws.Range("C18").Value = 0.80208333333333337
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = "6 h PostgreSQL-serverin kanssa säheltämistä. https://github.com/Eetu95/Open-source-IdM-solution/blob/master/Eetun%20muistiinpanoja/Ty%C3%B6t%20-%2010.10.2018.txt. 13 min ohjauskokous."

$ws.Rows.Item(18).AutoFit()

$ws.Range("F18").Select()
